$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style donor cells: A2 carries the "Indice" column style (bold + border),
# E2 carries the datetime number-format style used by data_partida.
$styleColA = $ws.Cells.Item(2,1)
$styleColE = $ws.Cells.Item(2,5)

# ---- Row 137 (Indice 136) ----
$styleColA.Copy()
$ws.Cells.Item(137,1).PasteSpecial(-4122)
$styleColE.Copy()
$ws.Cells.Item(137,5).PasteSpecial(-4122)
$ws.Cells.Item(137,1).Value = 136
$ws.Cells.Item(137,2).Value = "poland"
$ws.Cells.Item(137,3).Value = "ekstraklasa"
$ws.Cells.Item(137,4).Value = "2023-2024"
$ws.Cells.Item(137,5).Value = 45255.625
$ws.Cells.Item(137,6).Value = "Rakow"
$ws.Cells.Item(137,7).Value = 1
$ws.Cells.Item(137,8).Value = "Cracovia"
$ws.Cells.Item(137,9).Value = 1
$ws.Cells.Item(137,10).Value = 1.51
$ws.Cells.Item(137,11).Value = "18/11/2023 14:12"
$ws.Cells.Item(137,12).Value = 1.6
$ws.Cells.Item(137,13).Value = "25/11/2023 14:57"
$ws.Cells.Item(137,14).Value = 4.28
$ws.Cells.Item(137,15).Value = "18/11/2023 14:12"
$ws.Cells.Item(137,16).Value = 4.06
$ws.Cells.Item(137,17).Value = "25/11/2023 14:57"
$ws.Cells.Item(137,18).Value = 5.84
$ws.Cells.Item(137,19).Value = "18/11/2023 14:12"
$ws.Cells.Item(137,20).Value = 6.02
$ws.Cells.Item(137,21).Value = "25/11/2023 14:57"
$ws.Cells.Item(137,22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/rakow-czestochowa-cracovia/fgpCsgxj/"

# ---- Row 138 (Indice 137) ----
$styleColA.Copy()
$ws.Cells.Item(138,1).PasteSpecial(-4122)
$styleColE.Copy()
$ws.Cells.Item(138,5).PasteSpecial(-4122)
$ws.Cells.Item(138,1).Value = 137
$ws.Cells.Item(138,2).Value = "poland"
$ws.Cells.Item(138,3).Value = "ekstraklasa"
$ws.Cells.Item(138,4).Value = "2023-2024"
$ws.Cells.Item(138,5).Value = 45255.72916666666
$ws.Cells.Item(138,6).Value = "Pogon Szczecin"
$ws.Cells.Item(138,7).Value = 2
$ws.Cells.Item(138,8).Value = "Stal Mielec"
$ws.Cells.Item(138,9).Value = 3
$ws.Cells.Item(138,10).Value = 1.37
$ws.Cells.Item(138,11).Value = "17/11/2023 16:42"
$ws.Cells.Item(138,12).Value = 1.28
$ws.Cells.Item(138,13).Value = "25/11/2023 16:52"
$ws.Cells.Item(138,14).Value = 5.03
$ws.Cells.Item(138,15).Value = "17/11/2023 16:42"
$ws.Cells.Item(138,16).Value = 6.06
$ws.Cells.Item(138,17).Value = "25/11/2023 17:01"
$ws.Cells.Item(138,18).Value = 8.27
$ws.Cells.Item(138,19).Value = "17/11/2023 16:42"
$ws.Cells.Item(138,20).Value = 10.56
$ws.Cells.Item(138,21).Value = "25/11/2023 17:00"
$ws.Cells.Item(138,22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/pogon-szczecin-stal-mielec/GpfxleaM/"

# ---- Row 139 (Indice 138) ----
$styleColA.Copy()
$ws.Cells.Item(139,1).PasteSpecial(-4122)
$styleColE.Copy()
$ws.Cells.Item(139,5).PasteSpecial(-4122)
$ws.Cells.Item(139,1).Value = 138
$ws.Cells.Item(139,2).Value = "poland"
$ws.Cells.Item(139,3).Value = "ekstraklasa"
$ws.Cells.Item(139,4).Value = "2023-2024"
$ws.Cells.Item(139,5).Value = 45255.83333333334
$ws.Cells.Item(139,6).Value = "Legia"
$ws.Cells.Item(139,7).Value = 2
$ws.Cells.Item(139,8).Value = "Warta Poznan"
$ws.Cells.Item(139,9).Value = 2
$ws.Cells.Item(139,10).Value = 1.55
$ws.Cells.Item(139,11).Value = "18/11/2023 19:12"
$ws.Cells.Item(139,12).Value = 1.45
$ws.Cells.Item(139,13).Value = "25/11/2023 19:51"
$ws.Cells.Item(139,14).Value = 4.03
$ws.Cells.Item(139,15).Value = "18/11/2023 19:12"
$ws.Cells.Item(139,16).Value = 4.27
$ws.Cells.Item(139,17).Value = "25/11/2023 19:53"
$ws.Cells.Item(139,18).Value = 6.56
$ws.Cells.Item(139,19).Value = "18/11/2023 19:12"
$ws.Cells.Item(139,20).Value = 8.66
$ws.Cells.Item(139,21).Value = "25/11/2023 19:53"
$ws.Cells.Item(139,22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/legia-warta-poznan/xfeYlypG/"

# ---- Row 140 (Indice 139) ----
$styleColA.Copy()
$ws.Cells.Item(140,1).PasteSpecial(-4122)
$styleColE.Copy()
$ws.Cells.Item(140,5).PasteSpecial(-4122)
$ws.Cells.Item(140,1).Value = 139
$ws.Cells.Item(140,2).Value = "poland"
$ws.Cells.Item(140,3).Value = "ekstraklasa"
$ws.Cells.Item(140,4).Value = "2023-2024"
$ws.Cells.Item(140,5).Value = 45256.52083333334
$ws.Cells.Item(140,6).Value = "Puszcza"
$ws.Cells.Item(140,7).Value = 2
$ws.Cells.Item(140,8).Value = "Gornik Zabrze"
$ws.Cells.Item(140,9).Value = 1
$ws.Cells.Item(140,10).Value = 3.13
$ws.Cells.Item(140,11).Value = "19/11/2023 11:43"
$ws.Cells.Item(140,12).Value = 2.99
$ws.Cells.Item(140,13).Value = "26/11/2023 12:29"
$ws.Cells.Item(140,14).Value = 3.33
$ws.Cells.Item(140,15).Value = "19/11/2023 11:43"
$ws.Cells.Item(140,16).Value = 3.28
$ws.Cells.Item(140,17).Value = "26/11/2023 12:29"
$ws.Cells.Item(140,18).Value = 2.27
$ws.Cells.Item(140,19).Value = "19/11/2023 11:43"
$ws.Cells.Item(140,20).Value = 2.54
$ws.Cells.Item(140,21).Value = "26/11/2023 12:29"
$ws.Cells.Item(140,22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/puszcza-gornik-zabrze/AL8umFES/"

# ---- Row 141 (Indice 140) ----
$styleColA.Copy()
$ws.Cells.Item(141,1).PasteSpecial(-4122)
$styleColE.Copy()
$ws.Cells.Item(141,5).PasteSpecial(-4122)
$ws.Cells.Item(141,1).Value = 140
$ws.Cells.Item(141,2).Value = "poland"
$ws.Cells.Item(141,3).Value = "ekstraklasa"
$ws.Cells.Item(141,4).Value = "2023-2024"
$ws.Cells.Item(141,5).Value = 45256.625
$ws.Cells.Item(141,6).Value = "Ruch Chorzow"
$ws.Cells.Item(141,7).Value = 1
$ws.Cells.Item(141,8).Value = "Korona Kielce"
$ws.Cells.Item(141,9).Value = 1
$ws.Cells.Item(141,10).Value = 2.47
$ws.Cells.Item(141,11).Value = "19/11/2023 14:12"
$ws.Cells.Item(141,12).Value = 3.08
$ws.Cells.Item(141,13).Value = "26/11/2023 14:59"
$ws.Cells.Item(141,14).Value = 3.47
$ws.Cells.Item(141,15).Value = "19/11/2023 14:12"
$ws.Cells.Item(141,16).Value = 3.44
$ws.Cells.Item(141,17).Value = "26/11/2023 14:59"
$ws.Cells.Item(141,18).Value = 2.73
$ws.Cells.Item(141,19).Value = "19/11/2023 14:12"
$ws.Cells.Item(141,20).Value = 2.4
$ws.Cells.Item(141,21).Value = "26/11/2023 14:59"
$ws.Cells.Item(141,22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/ruch-chorzow-korona-kielce/8dlGtDid/"

# ---- Row 142 (Indice 141) ----
$styleColA.Copy()
$ws.Cells.Item(142,1).PasteSpecial(-4122)
$styleColE.Copy()
$ws.Cells.Item(142,5).PasteSpecial(-4122)
$ws.Cells.Item(142,1).Value = 141
$ws.Cells.Item(142,2).Value = "poland"
$ws.Cells.Item(142,3).Value = "ekstraklasa"
$ws.Cells.Item(142,4).Value = "2023-2024"
$ws.Cells.Item(142,5).Value = 45256.72916666666
$ws.Cells.Item(142,6).Value = "Lech Poznan"
$ws.Cells.Item(142,7).Value = 1
$ws.Cells.Item(142,8).Value = "Widzew Lodz"
$ws.Cells.Item(142,9).Value = 3
$ws.Cells.Item(142,10).Value = 1.44
$ws.Cells.Item(142,11).Value = "18/11/2023 18:13"
$ws.Cells.Item(142,12).Value = 1.55
$ws.Cells.Item(142,13).Value = "26/11/2023 17:29"
$ws.Cells.Item(142,14).Value = 4.5
$ws.Cells.Item(142,15).Value = "18/11/2023 18:13"
$ws.Cells.Item(142,16).Value = 4.28
$ws.Cells.Item(142,17).Value = "26/11/2023 17:29"
$ws.Cells.Item(142,18).Value = 6.44
$ws.Cells.Item(142,19).Value = "18/11/2023 18:13"
$ws.Cells.Item(142,20).Value = 6.34
$ws.Cells.Item(142,21).Value = "26/11/2023 17:29"
$ws.Cells.Item(142,22).Value = "https://www.betexplorer.com/football/poland/ekstraklasa/lech-poznan-widzew-lodz/fV5TkHU9/"

$excel.CutCopyMode = 0
Write-Output "Added rows 137-142"
